$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 1343
$ws.Range("I55").Value = 448.33334
$ws.Range("J55").Value = 1879.8
$ws.Range("K55").Value = 448.33334
$ws.Range("L55").Value = 1879.8
$ws.Range("M55").Value = -234.33334
$ws.Range("N55").Value = -2307.8
# Row 129
$ws.Range("H129").Value = 2190.8064
$ws.Range("I129").Value = 2814.6667
$ws.Range("J129").Value = 2041.08
$ws.Range("K129").Value = 8444.000100000001
$ws.Range("L129").Value = 6123.24
$ws.Range("M129").Value = -3444.000100000001
$ws.Range("N129").Value = -16123.24
# Row 137
$ws.Range("H137").Value = 4013.6592
$ws.Range("I137").Value = 1076.4348
$ws.Range("J137").Value = 7230.619
$ws.Range("K137").Value = 3229.3044
$ws.Range("L137").Value = 21691.857
$ws.Range("M137").Value = -679.3044
$ws.Range("N137").Value = -26791.857
# Row 138
$ws.Range("H138").Value = 1995.1039
$ws.Range("I138").Value = 2354.111
$ws.Range("K138").Value = 7062.333
$ws.Range("M138").Value = -1922.333

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 109
$ws.Range("H109").Value = 46784.332
$ws.Range("J109").Value = 46784.332
$ws.Range("L109").Value = 46784.332
$ws.Range("N109").Value = -49558.332
# Row 113
$ws.Range("H113").Value = 43843.5
$ws.Range("J113").Value = 43843.5
$ws.Range("L113").Value = 43843.5
$ws.Range("N113").Value = -52521.5
# Row 114
$ws.Range("H114").Value = 41931
$ws.Range("J114").Value = 41931
$ws.Range("L114").Value = 41931
$ws.Range("N114").Value = -50609
# Row 117
$ws.Range("H117").Value = 46350.4
$ws.Range("J117").Value = 46350.4
$ws.Range("L117").Value = 46350.4
$ws.Range("N117").Value = -55528.4
# Row 122
$ws.Range("H122").Value = 2862.4
$ws.Range("I122").Value = 2924
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 8772
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -6322
$ws.Range("N122").Value = -10900
# Row 132
$ws.Range("H132").Value = 13890224
$ws.Range("I132").Value = 21740090
$ws.Range("J132").Value = 1997.8462
$ws.Range("K132").Value = 65220270
$ws.Range("L132").Value = 5993.5386
$ws.Range("M132").Value = -65217740
$ws.Range("N132").Value = -11053.5386
# Row 140
$ws.Range("H140").Value = 31666
$ws.Range("J140").Value = 31666
$ws.Range("L140").Value = 31666
$ws.Range("N140").Value = -42026

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 184425.69
$ws.Range("I31").Value = 2099.4375
$ws.Range("J31").Value = 231477.62
$ws.Range("K31").Value = 2099.4375
$ws.Range("L31").Value = 231477.62
$ws.Range("M31").Value = -1804.4375
$ws.Range("N31").Value = -232067.62
# Row 34
$ws.Range("H34").Value = 184425.69
$ws.Range("I34").Value = 2099.4375
$ws.Range("J34").Value = 231477.62
$ws.Range("K34").Value = 2099.4375
$ws.Range("L34").Value = 231477.62
$ws.Range("M34").Value = -1897.4375
$ws.Range("N34").Value = -231881.62
# Row 86
$ws.Range("H86").Value = 8000
$ws.Range("I86").Value = 8000
$ws.Range("K86").Value = 8000
$ws.Range("M86").Value = -6877
# Row 89
$ws.Range("H89").Value = 8000
$ws.Range("I89").Value = 8000
$ws.Range("K89").Value = 40000
$ws.Range("M89").Value = -34384
# Row 112
$ws.Range("H112").Value = 32017.8
$ws.Range("J112").Value = 32017.8
$ws.Range("L112").Value = 32017.8
$ws.Range("N112").Value = -34971.8

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 2801.9688
$ws.Range("I131").Value = 8070.077
$ws.Range("J131").Value = 1459.1177
$ws.Range("K131").Value = 24210.231
$ws.Range("L131").Value = 4377.3531
$ws.Range("M131").Value = -19170.231
$ws.Range("N131").Value = -14457.3531

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 1633.3334
$ws.Range("I122").Value = 1487.5
$ws.Range("J122").Value = 2800
$ws.Range("K122").Value = 4462.5
$ws.Range("L122").Value = 8400
$ws.Range("M122").Value = -2012.5
$ws.Range("N122").Value = -13300
# Row 123
$ws.Range("H123").Value = 19882.666
$ws.Range("J123").Value = 19882.666
$ws.Range("L123").Value = 19882.666
$ws.Range("N123").Value = -24782.666
# Row 141
$ws.Range("H141").Value = 31703.23
$ws.Range("J141").Value = 31703.23
$ws.Range("L141").Value = 31703.23
$ws.Range("N141").Value = -42063.23

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 41
$ws.Range("H41").Value = 8000
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 10000
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 10000
$ws.Range("M41").Value = -1562
$ws.Range("N41").Value = -10876
# Row 68
$ws.Range("H68").Value = 3399.8
$ws.Range("I68").Value = 2999.5
$ws.Range("J68").Value = 3666.6667
$ws.Range("K68").Value = 2999.5
$ws.Range("L68").Value = 3666.6667
$ws.Range("M68").Value = -2250.5
$ws.Range("N68").Value = -5164.6667
# Row 71
$ws.Range("H71").Value = 3399.8
$ws.Range("I71").Value = 2999.5
$ws.Range("J71").Value = 3666.6667
$ws.Range("K71").Value = 14997.5
$ws.Range("L71").Value = 18333.3335
$ws.Range("M71").Value = -11253.5
$ws.Range("N71").Value = -25821.3335
# Row 112
$ws.Range("H112").Value = 37113.25
$ws.Range("J112").Value = 37113.25
$ws.Range("L112").Value = 37113.25
$ws.Range("N112").Value = -40067.25
# Row 119
$ws.Range("H119").Value = 47408
$ws.Range("J119").Value = 47408
$ws.Range("L119").Value = 47408
$ws.Range("N119").Value = -57084
# Row 132
$ws.Range("H132").Value = 3032.1052
$ws.Range("I132").Value = 1196.5238
$ws.Range("K132").Value = 3589.5714
$ws.Range("M132").Value = -1059.5714

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 110
$ws.Range("H110").Value = 22844.666
$ws.Range("J110").Value = 22844.666
$ws.Range("L110").Value = 22844.666
$ws.Range("N110").Value = -31024.666
# Row 112
$ws.Range("H112").Value = 36195.2
$ws.Range("J112").Value = 36195.2
$ws.Range("L112").Value = 36195.2
$ws.Range("N112").Value = -39149.2
# Row 119
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
# Row 122
$ws.Range("H122").Value = 9526143
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900
# Row 141
$ws.Range("H141").Value = 19254
$ws.Range("J141").Value = 19254
$ws.Range("L141").Value = 19254
$ws.Range("N141").Value = -29614
